# Re-shape the factors dataframe: insert 7 new factor rows (one-hot encoded
# close_code_* columns) and reorder/update existing rows so X/y can be
# derived from the sheet for modelling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 11-21 down by 7 rows (-> rows 18-28) to make room
# for the new close_code_* one-hot factor rows.
$ws.Rows("11:17").Insert()

# Final contents for rows 11-27 (factor, variable_type, dtype,
# unique_values, chi, p). Row 28 (user_dissatisfied) already landed in the
# right place after the insert above and keeps its original values.
$data = @(
    @("close_code_No Resolution Action", "analyse", "uint8", 2, [double]"212.5657420134024", [double]"3.785379350270401e-48"),
    @("assignment_group_company", "analyse2", "object", 10, [double]"108.405148086512", [double]"3.104338295130162e-19"),
    @("close_code_Data Correction", "analyse", "uint8", 2, [double]"103.9629160967067", [double]"2.061384633699008e-24"),
    @("priority_is_4", "analyse", "int64", 2, [double]"103.8263314178949", [double]"2.208504274319012e-24"),
    @("caller_is_employee", "analyse", "int64", 2, [double]"89.71333796111035", [double]"2.752927944834831e-21"),
    @("ka_count_log", "ignore", "int64", 9, [double]"69.71864170960302", [double]"5.590814295923098e-12"),
    @("contact_type", "ignore", "object", 5, [double]"62.83893362299189", [double]"7.336768356839108e-13"),
    @("breached_reason_code", "ignore", "object", 13, [double]"55.22842804766915", [double]"1.647005452924161e-07"),
    @("incident_has_ka_related_flag", "analyse", "int64", 2, [double]"38.95749777616205", [double]"4.33133544018824e-10"),
    @("self_service", "analyse", "int64", 2, [double]"31.27813668428506", [double]"2.23583427847182e-08"),
    @("close_code_Reboot / Restart", "analyse", "uint8", 2, [double]"29.17596056897717", [double]"6.609418315501527e-08"),
    @("appl_tier", "ignore", "object", 4, [double]"27.66288814486726", [double]"4.274586770349934e-06"),
    @("close_code_Security Modification", "analyse", "uint8", 2, [double]"21.3501624171497", [double]"3.825883411918223e-06"),
    @("close_code_Software Correction", "analyse", "uint8", 2, [double]"12.15795061766993", [double]"0.0004887905515663779"),
    @("close_code_Environmental Restoration", "analyse", "uint8", 2, [double]"3.50656659638963", [double]"0.06112600893982762"),
    @("close_code_Information Provided / Training", "analyse", "uint8", 2, [double]"0.3786840534703145", [double]"0.538308351370397"),
    @("caller_vip", "ignore", "int64", 2, [double]"0.2128688205425569", [double]"0.644528084295426")
)

$r = 11
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

Write-Host "Factors sheet reshaped: rows now $($ws.UsedRange.Rows.Count)"
